$d = $word.ActiveDocument

# 1) "Title: Prototyping Labs " -> "Title: Prototyping Lab " (remove the "s")
$d.Content.Find.Execute("Title: Prototyping Labs ", $true, $false, $false, $false, $false, $true, 1, $false, "Title: Prototyping Lab ", 2) | Out-Null

# 2) "(if applicable)" run merge / proofErr cleanup - replace the phrase with itself to normalize runs
$d.Content.Find.Execute("(if applicable)", $true, $false, $false, $false, $false, $true, 1, $false, "(if applicable)", 2) | Out-Null

# 3) Header: "Prototyping Labs at GIX" -> "Prototyping Lab at GIX"
$d.Content.Find.Execute("Prototyping Labs at GIX", $true, $false, $false, $false, $false, $true, 1, $false, "Prototyping Lab at GIX", 2) | Out-Null
